$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "ACORN Participating Countries"
$ws.Range("B11").Value = "Các nước tham gia vào nghiên cứu ACORN"
$ws.Range("A12").Value = "All 'orgname' are provided."
$ws.Range("B12").Value = "Tất cả `"orgname`" được cung cấp"
$ws.Range("A13").Value = "All 'patid' are provided."
$ws.Range("B13").Value = "Tất cả `"patid`" được cung cấp"
$ws.Range("A14").Value = "All 'specdate' are provided."
$ws.Range("B14").Value = "Tất cả `"specdate`" được cung cấp"
$ws.Range("A15").Value = "All 'specdate' are today or before today."
$ws.Range("B15").Value = "Tât cả `"specdate`" là ngày hôm này hoặc trước ngày hôm nay"
$ws.Range("A16").Value = "All 'specgroup' are provided."
$ws.Range("B16").Value = "Tất cả `"specgroup`" được cung cấp"
$ws.Range("A17").Value = "All 'specid' are provided."
$ws.Range("B17").Value = "Tất cả `"specid`" được cung cấp"
$ws.Range("A18").Value = "All dates of enrolment for HAI patients have a matching date in the HAI survey dataset"
$ws.Range("B18").Value = "Tất cả các ngày thu tuyển của bệnh nhân HAI phải tương ứng với ngày thực hiện giám sát điểm đã được định sẵn"
$ws.Range("A19").Value = "All Other Organisms"
$ws.Range("B19").Value = "Tất cả các vi sinh vật khác"
$ws.Range("A20").Value = "All valid records have an ACORN ID."
$ws.Range("B20").Value = "TBT"
$ws.Range("A21").Value = "AMR"
$ws.Range("B21").Value = "TBT"
$ws.Range("A22").Value = "and generate enrolment log."
$ws.Range("B22").Value = "và tạo nhật ký thu tuyển."
$ws.Range("A23").Value = "Attempting to connect."
$ws.Range("B23").Value = "Đang cố gắng kết nối."
$ws.Range("A24").Value = "Blood culture collected within 24 hours of admission (CAI) / symptom onset (HAI)"
$ws.Range("B24").Value = "Mẫu cấy máu thu thập trong vòng 24h sau khi nhập viện (CAI)/ khởi phát triệu chứng (HAI)"
$ws.Range("A25").Value = "Blood Culture Contaminants"
$ws.Range("B25").Value = "Tạp nhiễm cấy máu"
$ws.Range("A26").Value = "Bloodstream Infection (BSI)"
$ws.Range("B26").Value = "Nhiễm trùng huyết (BSI)"
$ws.Range("A27").Value = "Calculated age is consistent with 'Age Category'"
$ws.Range("B27").Value = "Tuổi được tính toán nhất quán với `"Phân loại tuổi`""
$ws.Range("A28").Value = "Calculated age isn't always consistent with 'Age Category'"
$ws.Range("B28").Value = "Tuổi được tính toán không nhất quán với `"Phân loại tuổi`""
$ws.Range("A29").Value = "Cancel"
$ws.Range("B29").Value = "Hủy bỏ"
$ws.Range("A30").Value = "Care should be taken when interpreting rates and AMR profiles where there are small numbers of cases or bacterial isolates: point estimates may be unreliable."
$ws.Range("B30").Value = "Cần thận trọng khi phiên giải về tỷ lệ và thông tin AMR khi có một số lượng nhỏ ca bệnh hoặc vi khuẩn được phân lập: ước tính điểm có thể không đáng tin cậy."
$ws.Range("A31").Value = "Clinical and day-28 outcomes are consistent."
$ws.Range("B31").Value = "Kêt quả lâm sàng và ngày 28 là nhất quán"
$ws.Range("A32").Value = "Clinical and day-28 outcomes aren't consistent for some dead patients."
$ws.Range("B32").Value = "Kết quả lâm sàng và ngày 28 là không nhất quán đối với một số bệnh nhân tử vong"
$ws.Range("A33").Value = "Clinical Outcome"
$ws.Range("B33").Value = "Kết cục lâm sàng"
$ws.Range("A34").Value = "Clinical Outcome Status:"
$ws.Range("B34").Value = "Tình trạng kết cục lâm sàng"
$ws.Range("A35").Value = "Co-resistances"
$ws.Range("B35").Value = "TBT"
$ws.Range("A36").Value = "Combine Susceptible + Intermediate"
$ws.Range("B36").Value = "Kết hợp Nhạy cảm + Trung bình"
$ws.Range("A37").Value = "Consider saving .acorn file on the cloud for additional security."
$ws.Range("B37").Value = "Cân nhắc lưu tệp .acorn trên đám mây vì mục đích bảo mật."
$ws.Range("A38").Value = "Contains names of organisms before and after mapping."
$ws.Range("B38").Value = "TBT"
$ws.Range("A39").Value = "Couldn't connect to server. Please check internet access."
$ws.Range("B39").Value = "Không thể kết nối với máy chủ. Vui lòng kiểm tra kết nối internet"
$ws.Range("A40").Value = "Critical errors with clinical data."
$ws.Range("B40").Value = "Lỗi nghiêm trọng với dữ liệu lâm sàng."
$ws.Range("A41").Value = "Culture results per specimen type"
$ws.Range("B41").Value = "Kết quả nuôi cấy cho mỗi loại mẫu"
$ws.Range("A42").Value = "Data Management"
$ws.Range("B42").Value = "Quản lý dữ liệu"
$ws.Range("A43").Value = "Date of Enrolment"
$ws.Range("B43").Value = "Ngày thu tuyển"
$ws.Range("A44").Value = "Day 28"
$ws.Range("B44").Value = "Ngày 28"
$ws.Range("A45").Value = "Day 28 Status:"
$ws.Range("B45").Value = "Tình trạng ngày 28"
$ws.Range("A46").Value = "Diagnosis at Enrolment"
$ws.Range("B46").Value = "Chẩn đoán tại thời điểm thu tuyển"
$ws.Range("A47").Value = "Dismiss"
$ws.Range("B47").Value = "TBT"
$ws.Range("A48").Value = "Distribution of Enrolments"
$ws.Range("B48").Value = "Phân bố thu tuyển"
$ws.Range("A49").Value = "Download Enrolment Log (.xlsx)"
$ws.Range("B49").Value = "Tải xuống Sổ thu tuyển (.xlsx)"
$ws.Range("A50").Value = "Download Lab Log (.xlsx)"
$ws.Range("B50").Value = "TBT"
$ws.Range("A70").Value = "HAI point prevalence by "
$ws.Range("B70").Value = "TBT"
$ws.Range("A110").Value = "Remove 'Not Cultured' specimens"
$ws.Range("B110").Value = "TBT"
$ws.Range("A111").Value = "Remove blood culture contaminants from the following visualizations"
$ws.Range("B111").Value = "Loại bỏ các tạp nhiễm trong quá trình cấy máu khỏi các hình ảnh trực quan"
$ws.Range("A112").Value = "Reset Enrolments Filters"
$ws.Range("B112").Value = "Cài đặt lại bộ lọc thu tuyển"
$ws.Range("A113").Value = "Resistance to 3rd gen. Cephalosporins Over Time"
$ws.Range("B113").Value = "Kháng Cephalosporins thế hệ 3 theo thời gian."
$ws.Range("A114").Value = "Resistance to Carbapenems Over Time"
$ws.Range("B114").Value = "Kháng Carbapenems theo thời gian"
$ws.Range("A115").Value = "Resistance to Fluoroquinolones Over Time"
$ws.Range("B115").Value = "Kháng Fluoroquinolones theo thời gian"
$ws.Range("A116").Value = "Resistance to Oxacillin Over Time"
$ws.Range("B116").Value = "Kháng Oxacillin theo thời gian"
$ws.Range("A117").Value = "Resistance to Penicillin G - meningitis Over Time"
$ws.Range("B117").Value = "Kháng Penicillin G -meningitis theo thời gian"
$ws.Range("A118").Value = "Resistance to Penicillin G Over Time"
$ws.Range("B118").Value = "Kháng Penicillin G theo thời gian"
$ws.Range("A119").Value = "Retriving data from REDCap server."
$ws.Range("B119").Value = "Truy xuất dữ liệu từ máy chủ REDCap."
$ws.Range("A120").Value = "Save .acorn file"
$ws.Range("B120").Value = "Lưu tệp .acorn"
$ws.Range("A121").Value = "Save acorn data"
$ws.Range("B121").Value = "Lưu dữ liệu acorn"
$ws.Range("A122").Value = "Save on Server"
$ws.Range("B122").Value = "Lưu trên máy chủ"
$ws.Range("A123").Value = "See Breakdown by Ward"
$ws.Range("B123").Value = "Xem Phân tích theo Khoa"
$ws.Range("A124").Value = "See by Week"
$ws.Range("B124").Value = "Xem theo Tuần"
